# Reset commit counts to 0 for all group members/leaders that had a
# non-zero count, keeping the "<name> : <count>" text format intact.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "LeeYinWei : 0"
$ws.Range("C2").Value = "unknown899 : 0"
$ws.Range("B3").Value = "neoAurora : 0"
$ws.Range("D3").Value = "howardhung14 : 0"
$ws.Range("B4").Value = "yoyo0213 : 0"
$ws.Range("C4").Value = "JonathanYangSW : 0"
$ws.Range("D4").Value = "GinoChen113511247 : 0"
$ws.Range("C5").Value = "peienwu1216 : 0"
$ws.Range("D5").Value = "chxyuuu : 0"
$ws.Range("B6").Value = "ginny923 : 0"
$ws.Range("C6").Value = "joanna0420 : 0"
$ws.Range("D6").Value = "dua0505 : 0"
$ws.Range("B7").Value = "jui-pixel : 0"
$ws.Range("D7").Value = "charles691 : 0"
$ws.Range("B8").Value = "Tony104147 : 0"
$ws.Range("B9").Value = "haleychang0530 : 0"
$ws.Range("C9").Value = "Hazel-1212 : 0"
$ws.Range("B10").Value = "CHENG-JE : 0"
$ws.Range("C10").Value = "lwc-ed : 0"
$ws.Range("B11").Value = "tpvupu : 0"
$ws.Range("C11").Value = "xiaotin22 : 0"
$ws.Range("D11").Value = "calistayang : 0"
$ws.Range("B13").Value = "kufanghua : 0"
$ws.Range("C13").Value = "yezh0915 : 0"
$ws.Range("D13").Value = "fiesta0217 : 0"
$ws.Range("D14").Value = "jing1688 : 0"
$ws.Range("B15").Value = "weiouo-0817 : 0"
$ws.Range("B16").Value = "gamemode0701 : 0"
$ws.Range("C16").Value = "Tonyyu2403 : 0"
$ws.Range("B17").Value = "TerryCheese : 0"
$ws.Range("C17").Value = "junlin27 : 0"
$ws.Range("B18").Value = "Miiaow3011 : 0"
$ws.Range("C18").Value = "bonnieliao774 : 0"
$ws.Range("D18").Value = "emmazheng0318 : 0"
$ws.Range("C19").Value = "TedChueh : 0"
$ws.Range("B20").Value = "max052028 : 0"
$ws.Range("B21").Value = "houyuankai : 0"
$ws.Range("B22").Value = "0u88 : 0"
